$d = $word.ActiveDocument

$replacements = @(
    @("2023-08-19 Saturday", "2023-08-20 Sunday"),
    @("86÷7=12, 2", "72÷3=24, 0"),
    @("10÷7=1, 3", "20÷6=3, 2"),
    @("90÷3=30, 0", "42÷5=8, 2"),
    @("56÷8=7, 0", "55÷9=6, 1"),
    @("74÷3=24, 2", "96÷4=24, 0"),
    @("17÷6=2, 5", "60÷7=8, 4"),
    @("59÷7=8, 3", "57÷5=11, 2"),
    @("56÷9=6, 2", "26÷6=4, 2"),
    @("62÷8=7, 6", "78÷8=9, 6"),
    @("99÷7=14, 1", "87÷2=43, 1"),
    @("83÷6=13, 5", "84÷9=9, 3"),
    @("40÷7=5, 5", "92÷4=23, 0"),
    @("44÷6=7, 2", "13÷5=2, 3"),
    @("79÷9=8, 7", "79÷4=19, 3"),
    @("75÷9=8, 3", "44÷9=4, 8"),
    @("71÷3=23, 2", "57÷9=6, 3"),
    @("77÷2=38, 1", "52÷7=7, 3"),
    @("72÷4=18, 0", "41÷7=5, 6"),
    @("89÷9=9, 8", "93÷5=18, 3"),
    @("59÷3=19, 2", "42÷8=5, 2"),
    @("33÷4=8, 1", "21÷7=3, 0"),
    @("29÷3=9, 2", "63÷9=7, 0"),
    @("59÷9=6, 5", "44÷5=8, 4"),
    @("25÷6=4, 1", "95÷7=13, 4")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
